$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Sector" (sheet1.xml): add column G = Sector_PT
# ---------------------------------------------------------------------------
$wsSector = $wb.Worksheets.Item("Sector")
$wsSector.Range("G1").Value = "Sector_PT"
$wsSector.Range("G2").Value = "Cuidados da Saúde"
$wsSector.Range("G3").Value = "Ciências Biológicas"

# ---------------------------------------------------------------------------
# Sheet "Industry" (sheet2.xml): add column G = Industry_PT
# ---------------------------------------------------------------------------
$wsIndustry = $wb.Worksheets.Item("Industry")
$wsIndustry.Range("G1").Value = "Industry_PT"
$wsIndustry.Range("G2").Value = "Cuidados da Saúde"
$wsIndustry.Range("G3").Value = "Dispositivos Médicos"
$wsIndustry.Range("G4").Value = "Produtos Farmacêuticos"

# ---------------------------------------------------------------------------
# Sheet "Product" (sheet3.xml): append pt rows (17-19)
# ---------------------------------------------------------------------------
$wsProduct = $wb.Worksheets.Item("Product")

$wsProduct.Range("A2:C2").Copy()
$wsProduct.Range("A17:C17").PasteSpecial(-4122)
$wsProduct.Range("A17").Value = 1
$wsProduct.Range("B17").Value = "pt"
$wsProduct.Range("C17").Value = "Produto A"

$wsProduct.Range("A3:C3").Copy()
$wsProduct.Range("A18:C18").PasteSpecial(-4122)
$wsProduct.Range("A18").Value = 2
$wsProduct.Range("B18").Value = "pt"
$wsProduct.Range("C18").Value = "Produto B"

$wsProduct.Range("A4:C4").Copy()
$wsProduct.Range("A19:C19").PasteSpecial(-4122)
$wsProduct.Range("A19").Value = 3
$wsProduct.Range("B19").Value = "pt"
$wsProduct.Range("C19").Value = "Produto C"

# ---------------------------------------------------------------------------
# Sheet "Company" (sheet4.xml): append pt rows (17-19)
# ---------------------------------------------------------------------------
$wsCompany = $wb.Worksheets.Item("Company")

$wsCompany.Range("A2:C2").Copy()
$wsCompany.Range("A17:C17").PasteSpecial(-4122)
$wsCompany.Range("A17").Value = 1
$wsCompany.Range("B17").Value = "pt"
$wsCompany.Range("C17").Value = "Empresa ABC"

$wsCompany.Range("A3:C3").Copy()
$wsCompany.Range("A18:C18").PasteSpecial(-4122)
$wsCompany.Range("A18").Value = 2
$wsCompany.Range("B18").Value = "pt"
$wsCompany.Range("C18").Value = "Empresa ACME"

$wsCompany.Range("A4:C4").Copy()
$wsCompany.Range("A19:C19").PasteSpecial(-4122)
$wsCompany.Range("A19").Value = 3
$wsCompany.Range("B19").Value = "pt"
$wsCompany.Range("C19").Value = "Empresa XYZ"

# ---------------------------------------------------------------------------
# Sheet "Language" (sheet6.xml): append pt row (7)
# ---------------------------------------------------------------------------
$wsLanguage = $wb.Worksheets.Item("Language")
$wsLanguage.Range("A7").Value = "pt"
$wsLanguage.Range("B7").Value = "Português"

# ---------------------------------------------------------------------------
# Sheet "UI_Strings" (sheet7.xml): insert new pt translation rows before the
# trailing themeRGB block (old rows 74-79 shift down to 87-92), then append
# one more themeRGB/pt row (93).
# ---------------------------------------------------------------------------
$wsUI = $wb.Worksheets.Item("UI_Strings")

$wsUI.Rows.Item(74).Resize(13).Insert()

$wsUI.Range("A74").Value = "Company"
$wsUI.Range("B74").Value = "Empresa"
$wsUI.Range("C74").Value = "pt"

$wsUI.Range("A75").Value = "Date"
$wsUI.Range("B75").Value = "Data"
$wsUI.Range("C75").Value = "pt"

$wsUI.Range("A76").Value = "Extended Price"
$wsUI.Range("B76").Value = "Preço Estendido"
$wsUI.Range("C76").Value = "pt"

$wsUI.Range("A77").Value = "Industry"
$wsUI.Range("B77").Value = "Industria"
$wsUI.Range("C77").Value = "pt"

$wsUI.Range("A78").Value = "Language"
$wsUI.Range("B78").Value = "Idioma"
$wsUI.Range("C78").Value = "pt"

$wsUI.Range("A79").Value = "Order #"
$wsUI.Range("B79").Value = "Nº Pedido"
$wsUI.Range("C79").Value = "pt"

$wsUI.Range("A80").Value = "Order Details"
$wsUI.Range("B80").Value = "Detalhe do Pedido"
$wsUI.Range("C80").Value = "pt"

$wsUI.Range("A81").Value = "Product"
$wsUI.Range("B81").Value = "Produto"
$wsUI.Range("C81").Value = "pt"

$wsUI.Range("A82").Value = "Quantity"
$wsUI.Range("B82").Value = "Quantidade"
$wsUI.Range("C82").Value = "pt"

$wsUI.Range("A83").Value = "Sector"
$wsUI.Range("B83").Value = "Setor"
$wsUI.Range("C83").Value = "pt"

$wsUI.Range("A84").Value = "Total Orders"
$wsUI.Range("B84").Value = "Total de Pedidos"
$wsUI.Range("C84").Value = "pt"

$wsUI.Range("A85").Value = "Unit Price"
$wsUI.Range("B85").Value = "Preço Unitário"
$wsUI.Range("C85").Value = "pt"

$wsUI.Range("A86").Value = " Select your preferred language… "
$wsUI.Range("A86").Value = "Select your preferred language…"
$wsUI.Range("B86").Value = " Selecione seu idioma de preferência…"
$wsUI.Range("C86").Value = "pt"

# New trailing themeRGB / pt row (93), copying the style of the last
# themeRGB row (now at 92, originally 79) which carries style s="4".
$wsUI.Range("B92").Copy()
$wsUI.Range("B93").PasteSpecial(-4122)
$wsUI.Range("A93").Value = "themeRGB"
$wsUI.Range("B93").Value = "RGB(48,125,48)"
$wsUI.Range("C93").Value = "pt"

# ---------------------------------------------------------------------------
# Selections / active sheet bookkeeping (matches the final sheetView state)
# ---------------------------------------------------------------------------
$wsSector.Range("G3").Select()
$wsIndustry.Range("G3").Select()
$wsCompany.Range("B26").Select()
$wsLanguage.Range("B8").Select()
$wsUI.Range("B85").Select()
$wsProduct.Range("C20").Select()
